# model.xlsx update: rename strategies "standard"/"new_treatment" to
# "seritinib"/"volantor", and add a new "cendralimab" strategy, including
# its per-strategy variable rows on the "variables" sheet.

$wb = $excel.ActiveWorkbook

# --- strategies sheet --------------------------------------------------
$strategies = $wb.Worksheets.Item("strategies")

$strategies.Range("A2").Value = "seritinib"
$strategies.Range("B2").Value = "Seritinib"
$strategies.Range("A3").Value = "volantor"
$strategies.Range("B3").Value = "Volantor"
$strategies.Range("A4").Value = "cendralimab"
$strategies.Range("B4").Value = "Cendralimab"

# --- variables sheet ----------------------------------------------------
$variables = $wb.Worksheets.Item("variables")

# update existing strategy-specific variable rows to use the new names
$variables.Range("D9").Value = "seritinib"
$variables.Range("D10").Value = "volantor"
$variables.Range("D11").Value = "seritinib"
$variables.Range("D12").Value = "volantor"

# Row 19: c_treatment for the new cendralimab strategy (clone row 9, which
# holds the other c_treatment/strategy row, then adjust the cells that
# differ).
$variables.Rows.Item(9).Copy()
$variables.Rows.Item(19).PasteSpecial()
$excel.CutCopyMode = $false

$variables.Range("D19").Value = "cendralimab"
$variables.Range("F19").Value = "gamma(mean = 5000, sd = 1000)"

# C19 needs to hold the *text* "5000" (matching how the other value cells
# on this sheet store their numbers as text). Assigning a numeric-looking
# string directly gets auto-converted to a real number, so instead we
# write a TEXT() formula and then flatten it back down to a plain value.
$variables.Range("C19").Formula = '=TEXT(5000,"0")'
$variables.Range("C19").Copy()
$variables.Range("C19").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row 20: hr_progression for the new cendralimab strategy (clone row 12,
# which holds the other hr_progression/strategy row with a sampling
# distribution already filled in).
$variables.Rows.Item(12).Copy()
$variables.Rows.Item(20).PasteSpecial()
$excel.CutCopyMode = $false

$variables.Range("D20").Value = "cendralimab"
$variables.Range("F20").Value = "lognormal(mean = 0.4, sd = 0.08)"

$variables.Range("C20").Formula = '=TEXT(0.4,"0.0")'
$variables.Range("C20").Copy()
$variables.Range("C20").PasteSpecial(-4163)
$excel.CutCopyMode = $false

Write-Host "strategies + variables updated"
